$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 4 days.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 4
}

# Update the "Actual Production (MW)" values in column B for rows 24-51
# (old portfolio data being removed / replaced).
$bValues = @{
    24 = 8
    25 = 22
    26 = 55
    27 = 91
    28 = 119
    29 = 150
    30 = 218
    31 = 283
    32 = 321
    33 = 355
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
}

foreach ($r in $bValues.Keys) {
    $ws.Cells.Item($r, 2).Value2 = $bValues[$r]
}
